$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Metadata sheet: bump the ConceptMap "Date" property value
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2026-01-23T08:28:04+00:00"

# ------------------------------------------------------------------
# 2. "Mapping Table 1" sheet: insert two new mapping rows
#    (route / site) right before the "consumable" row, pushing the
#    existing rows down by two. We avoid Rows.Insert() (it creates a
#    brand-new, unstyled cellXf) and instead grow the used range by
#    cloning the formatting of the last existing data row, then
#    rewrite rows 9..15 top-to-bottom. Rows 1-8 are untouched.
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Mapping Table 1")

# Grow the sheet by 2 rows, reusing the existing "data row" style (s=2)
# so no new cellXf gets appended to styles.xml.
$ws.Range("A13:E13").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("A15:E15").PasteSpecial(-4122)

# Final content for rows 9-15 (Source | Relationship | Target).
# Column B and E stay blank for these data rows.
$rows = @(
    @("FRCDAVaccinRecommande.routeCode", "equivalent", "FRImmunizationRecommendationDocument.supportingImmunization:FRImmunizationDocument.route"),
    @("FRCDAVaccinRecommande.approachSiteCode", "equivalent", "FRImmunizationRecommendationDocument.supportingImmunization:FRImmunizationDocument.site"),
    @("FRCDAVaccinRecommande.consumable.FRCDAProduitDeSante", "equivalent", "FRImmunizationRecommendationDocument.recommendation.vaccineCode"),
    @("FRCDAVaccinRecommande.consumable.FRCDAProduitDeSante.code.translation", "equivalent", "FRImmunizationRecommendationDocument.recommendation.vaccineCode.coding:translation"),
    @("FRCDAVaccinRecommande.entryRelationship:frPrescription", "equivalent", "FRImmunizationRecommendationDocument.recommendation.supportingPatientInformation"),
    @("FRCDAVaccinRecommande.entryRelationship:frRangDeLaVaccination", "equivalent", "FRImmunizationRecommendationDocument.recommendation.seriesDosesPositiveInt"),
    @("FRCDAVaccinRecommande.entryRelationship:frCommentaireER", "equivalent", "FRImmunizationRecommendationDocument.description")
)

$r = 9
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $r = $r + 1
}
